$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "Razon social"/"Nombre Fantasia" entries that used a comma as a
# separator between co-contractors; the comma should be a period.
$ws.Range("E2:F99").Replace(", ", ". ")

# Fix "Importe" amounts that were stored as text using the Argentine
# locale format (dot as thousands separator, comma as decimal
# separator). Convert them to a plain "1234.56" style numeric-looking
# text: drop the thousands separators, then turn the decimal comma
# into a decimal point.
$ws.Range("H2:H99").NumberFormat = "@"
$ws.Range("H2:H99").Replace(".", "")
$ws.Range("H2:H99").Replace(",", ".")
$ws.Range("H2:H99").Style = "Normal"
